$wb = $excel.ActiveWorkbook

# --- Rename the "Requested quantity" headers ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- Forecast data rows ---
$data = @(
    @(45508.99999999999, 12, 1.608169957894659, 22.01413993554768),
    @(45515.99999999999, 11, 0.4001372097370263, 21.50579877197319),
    @(45606.99999999999, 3, -6.660651181436743, 13.01052522826351),
    @(45613.99999999999, 2, -7.559472541081656, 11.84132175592505),
    @(45620.99999999999, 2, -8.485461640899413, 11.57959670483723),
    @(45627.99999999999, 1, -8.276463455834536, 11.73190867677988),
    @(45634.99999999999, 0, -9.675300265462408, 11.06329064881705),
    @(45641.99999999999, 0, -10.86390335892662, 9.206301378414606),
    @(45648.99999999999, 0, -11.00326408393548, 9.472954623114367),
    @(45655.99999999999, 0, -11.58865757059338, 7.858942678371361),
    @(45662.99999999999, 0, -11.54755500314883, 7.8286142368649)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# --- Match formatting used on the other sheets: bold/bordered/centered
# header style, and the date number format on column A ---
$wsWeekly.Range("B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A12").PasteSpecial(-4122)
